$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1:3").Delete()
